# Daily attendance processing - reorder "Recorded By" names.
#
# The "Recorded By" column (G) lists names as "System, <email>". For every
# session row whose recorder list is exactly "System, <single email>", the
# email should be moved to the front: "<email>, System". Rows whose second
# entry is the backup/service account (backup@backdoor.com) are left
# untouched, as are rows with a different shape (single name, or more than
# two names).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dim = $ws.UsedRange
$lastRow = $dim.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -eq $null) {
        continue
    }

    # Only touch simple two-party "System, <email>" values (skip the
    # three-party backup/system rows and the backdoor account).
    if ($val -like "System, *" -and $val -notlike "System, *, *" -and $val -notlike "*backup@backdoor.com*") {
        $other = $val.Substring(8)
        $newVal = "$other, System"
        $cell.Value = $newVal
    }
}
